# Generate Report for Handoff
# Updates the localization-status workbook to reflect that the report was
# (re)generated for a handoff rather than reflecting a handback:
#   - Status text changes from "Handed back: in sync with en-US" to
#     "Ready for handoff" on every sheet that shows it.
#   - The "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
#     timestamps are refreshed to the new generation time.
#   - The Status columns, now holding shorter text, are narrowed to match
#     Excel's auto-fit behavior.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-31 05:01:09"
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333336
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333336

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-31 05:00:59"
$wsZhCn.Columns.Item(3).ColumnWidth = 16.333333333333336

# --- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-31 05:01:09"
$wsDeDe.Columns.Item(3).ColumnWidth = 16.333333333333336
